$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "banned term" for row 2 (iphone 12 64gb) from "mini watch" to "mini watch 11"
$ws.Range("B2").Value = "mini watch 11"

# Update the selected cell to B3 (matches the selection change in the saved file)
$ws.Range("B3").Select()
